$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Change the "Runmode" column values in rows 3-7 from "Y" to "N"
# (row 2 keeps its existing "Y" value and is left untouched)
$ws.Range("C3:C7").Value = "N"

# Update the sheet's active cell / selection to match the edited range
$ws.Activate()
$ws.Range("C3:C7").Select()
